$wb = $excel.ActiveWorkbook

# --- Sheet1 (pcroprep) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D35").Value = [double]"4.4408920985006262E-16"
$ws1.Range("F35").Value = [double]"1.7763568394002505E-15"
$ws1.Range("G35").Value = -241.4

$ws1.Range("D39").Value = 940.06243687239612
$ws1.Range("F39").Value = 284.77792550326666
$ws1.Range("G39").Value = -521.82207449673342

# --- Sheet4 (pdietrep) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E6").Value = 1214.8943640250634
$ws4.Range("F6").Value = -942.89115018061943
$ws4.Range("G6").Value = 56.302832511705247

$ws4.Range("E7").Value = 41.227792490760066
$ws4.Range("F7").Value = -28.900398385590492
$ws4.Range("G7").Value = 58.789185883110214

$ws4.Range("E8").Value = 14.785465764643718
$ws4.Range("F8").Value = -49.948099661526754
$ws4.Range("G8").Value = 22.840493440001765

$ws4.Range("E9").Value = 215.13404923941326
$ws4.Range("F9").Value = -108.53377789143909
$ws4.Range("G9").Value = 66.467542092912097

# --- Sheet5 (pradar) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("D15").Value = [double]"1.7763568394002505E-15"
$ws5.Range("E15").Value = [double]"7.358561886496481E-16"
$ws5.Range("F15").Value = -241.4

# --- Sheet6 (plandrep) ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("S11").ClearContents()

# --- Sheet7 (plaborrep) ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("R3").ClearContents()
$ws7.Range("AF3").Value = 0.91221062026425459

# --- Sheet8 (pfertrep) ---
$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("S5").ClearContents()
$ws8.Range("Z5").Value = 291586.2928207317

$ws8.Range("S6").ClearContents()
$ws8.Range("Z6").Value = 341415.28265000007

$ws8.Range("S7").ClearContents()
$ws8.Range("Z7").Value = 323916.94579756097
